$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update regression-output table (columns B:H) for rows 2-29 with refreshed model coefficients.
# Rows whose Std.Err./t/P/CI columns (C:G) are blank in the source keep those cells blank;
# row 18 gains populated C:G values and rows 6 & 8 lose theirs, matching the refreshed fit.

# Row 2
$ws.Range("H2").Value = 0.1061941420915777
# Row 3
$ws.Range("B3").Value = 0.07538709310445521
$ws.Range("H3").Value = 0.1815812351960329
# Row 4
$ws.Range("B4").Value = 0.09129632098960513
$ws.Range("H4").Value = 0.1974904630811828
# Row 5
$ws.Range("B5").Value = 0.08244452408814788
$ws.Range("C5").Value = 0.009968173993340684
$ws.Range("D5").Value = 7.770984081288365
$ws.Range("E5").Value = 0.07090955337166575
$ws.Range("F5").Value = 0.06289128833426069
$ws.Range("G5").Value = 0.1019977598420342
$ws.Range("H5").Value = 0.1886386661797255
# Row 6
$ws.Range("B6").Value = 0.033075347815321
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 0.1392694899068987
# Row 7
$ws.Range("B7").Value = 0.02282775969911425
$ws.Range("C7").Value = 0.002877408574986257
$ws.Range("D7").Value = 1.413239325310184
$ws.Range("E7").Value = 0.05824660992423984
$ws.Range("F7").Value = 0.01718184450884557
$ws.Range("G7").Value = 0.02847367488938306
$ws.Range("H7").Value = 0.1290219017906919
# Row 8
$ws.Range("B8").Value = 0.01432766963596442
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0.1205218117275421
# Row 9
$ws.Range("B9").Value = 0.01128786288893247
$ws.Range("C9").Value = 0.006679985189627877
$ws.Range("D9").Value = -0.2362902133199316
$ws.Range("E9").Value = 0.05838131538722428
$ws.Range("F9").Value = -0.00189088772596585
$ws.Range("G9").Value = 0.02446661350383084
$ws.Range("H9").Value = 0.1174820049805101
# Row 10
$ws.Range("B10").Value = 0.005150575418067254
$ws.Range("C10").Value = 0.001907382030096153
$ws.Range("D10").Value = -0.5046112344498452
$ws.Range("E10").Value = 0.04019586879019919
$ws.Range("F10").Value = 0.001391177501802804
$ws.Range("G10").Value = 0.008909973334331656
$ws.Range("H10").Value = 0.1113447175096449
# Row 11
$ws.Range("B11").Value = 0.0238976330016012
$ws.Range("H11").Value = 0.1300917750931788
# Row 12
$ws.Range("B12").Value = 0.03790309372304096
$ws.Range("H12").Value = 0.1440972358146186
# Row 13
$ws.Range("B13").Value = 0.04458223464389077
$ws.Range("H13").Value = 0.1507763767354684
# Row 14
$ws.Range("B14").Value = 0.05086384210138638
$ws.Range("H14").Value = 0.157057984192964
# Row 15
$ws.Range("B15").Value = 0.05543922514084509
$ws.Range("H15").Value = 0.1616333672324227
# Row 16
$ws.Range("B16").Value = 0.05758888905792175
$ws.Range("C16").Value = 0.008283713963978647
$ws.Range("D16").Value = 11.85619081500656
$ws.Range("E16").Value = 0.05071761771923514
$ws.Range("F16").Value = 0.04133174130250004
$ws.Range("G16").Value = 0.07384603681334358
$ws.Range("H16").Value = 0.1637830311494994
# Row 17
$ws.Range("B17").Value = 0.06246454816749107
$ws.Range("C17").Value = 0.008043627690695396
$ws.Range("D17").Value = 12.54365822952609
$ws.Range("E17").Value = 0.04585919853021007
$ws.Range("F17").Value = 0.04668213259341834
$ws.Range("G17").Value = 0.07824696374156373
$ws.Range("H17").Value = 0.1686586902590687
# Row 18
$ws.Range("B18").Value = -0.1061941420915777
$ws.Range("C18").Value = 0.01228733386269351
$ws.Range("D18").Value = -18.60005276615585
$ws.Range("E18").Value = 0.03593855487942071
$ws.Range("F18").Value = -0.1303675710094518
$ws.Range("G18").Value = -0.08202071317370364
# Row 19
$ws.Range("B19").Value = 0.06571682889121765
$ws.Range("C19").Value = 0.0080271893419201
$ws.Range("D19").Value = 12.97719454180985
$ws.Range("E19").Value = 0.05245693064346862
$ws.Range("F19").Value = 0.04996430543512473
$ws.Range("G19").Value = 0.08146935234731062
$ws.Range("H19").Value = 0.1719109709827953
# Row 20
$ws.Range("B20").Value = 0.06642675937528832
$ws.Range("C20").Value = 0.008807015863611984
$ws.Range("D20").Value = 12.65377514829668
$ws.Range("E20").Value = 0.05440432688346138
$ws.Range("F20").Value = 0.04907456038323566
$ws.Range("G20").Value = 0.08377895836734112
$ws.Range("H20").Value = 0.172620901466866
# Row 21
$ws.Range("B21").Value = 0.06421143820889143
$ws.Range("C21").Value = 0.00850657440918221
$ws.Range("D21").Value = 12.43928540507902
$ws.Range("E21").Value = 0.05047108041171374
$ws.Range("F21").Value = 0.04750727662943188
$ws.Range("G21").Value = 0.08091559978835082
$ws.Range("H21").Value = 0.1704055803004691
# Row 22
$ws.Range("B22").Value = 0.06771490746122449
$ws.Range("C22").Value = 0.008362230008146013
$ws.Range("D22").Value = 12.6762303427851
$ws.Range("E22").Value = 0.05695710903125414
$ws.Range("F22").Value = 0.05129907030543083
$ws.Range("G22").Value = 0.08413074461701819
$ws.Range("H22").Value = 0.1739090495528021
# Row 23
$ws.Range("B23").Value = 0.0692506945484794
$ws.Range("C23").Value = 0.0083450524569745
$ws.Range("D23").Value = 13.14085127011487
$ws.Range("E23").Value = 0.05480444275960535
$ws.Range("F23").Value = 0.0528731132637595
$ws.Range("G23").Value = 0.08562827583319915
$ws.Range("H23").Value = 0.1754448366400571
# Row 24
$ws.Range("B24").Value = 0.07152841509902791
$ws.Range("C24").Value = 0.00776604083212342
$ws.Range("D24").Value = 13.35602331146976
$ws.Range("E24").Value = 0.05340495578331553
$ws.Range("F24").Value = 0.05628916459631818
$ws.Range("G24").Value = 0.08676766560173763
$ws.Range("H24").Value = 0.1777225571906056
# Row 25
$ws.Range("B25").Value = 0.0716482692395748
$ws.Range("C25").Value = 0.00831511932356521
$ws.Range("D25").Value = 12.69184647169015
$ws.Range("E25").Value = 0.06446798172904109
$ws.Range("F25").Value = 0.05532522009130807
$ws.Range("G25").Value = 0.08797131838784195
$ws.Range("H25").Value = 0.1778424113311525
# Row 26
$ws.Range("B26").Value = 0.0735420447311585
$ws.Range("C26").Value = 0.008276722301019701
$ws.Range("D26").Value = 12.81030019070673
$ws.Range("E26").Value = 0.05612003391098392
$ws.Range("F26").Value = 0.05729873759398849
$ws.Range("G26").Value = 0.08978535186832845
$ws.Range("H26").Value = 0.1797361868227361
# Row 27
$ws.Range("B27").Value = 0.07390342536206715
$ws.Range("C27").Value = 0.009049699157121269
$ws.Range("D27").Value = 12.54422792947498
$ws.Range("E27").Value = 0.05869854183001363
$ws.Range("F27").Value = 0.05613662179409668
$ws.Range("G27").Value = 0.09167022893003768
$ws.Range("H27").Value = 0.1800975674536448
# Row 28
$ws.Range("B28").Value = 0.0737856038310364
$ws.Range("C28").Value = 0.007757768211367496
$ws.Range("D28").Value = 12.66078995606631
$ws.Range("E28").Value = 0.0791336495691254
$ws.Range("F28").Value = 0.05856537247682368
$ws.Range("G28").Value = 0.08900583518524936
$ws.Range("H28").Value = 0.1799797459226141
# Row 29
$ws.Range("B29").Value = 0.02706431340998981
$ws.Range("C29").Value = 0.004886929033736671
$ws.Range("D29").Value = 0.8147582788376316
$ws.Range("E29").Value = 0.01181268684215629
$ws.Range("F29").Value = 0.01742668874187392
$ws.Range("G29").Value = 0.03670193807810559
$ws.Range("H29").Value = 0.1332584555015675

Write-Output "Updated cap_gen_year coefficient table (rows 2-29, columns B:H)."
